$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1999
$ws.Range("J48").Value = 1999
$ws.Range("L48").Value = 5997
$ws.Range("N48").Value = -6581
$ws.Range("H56").Value = 1999
$ws.Range("J56").Value = 1999
$ws.Range("L56").Value = 5997
$ws.Range("N56").Value = -7065
$ws.Range("H68").Value = 32295
$ws.Range("J68").Value = 32295
$ws.Range("L68").Value = 32295
$ws.Range("N68").Value = -33793
$ws.Range("H71").Value = 32295
$ws.Range("J71").Value = 32295
$ws.Range("L71").Value = 96885
$ws.Range("N71").Value = -104373
$ws.Range("H100").Value = 1695.8
$ws.Range("I100").Value = 1671.1765
$ws.Range("K100").Value = 1671.1765
$ws.Range("M100").Value = -1130.1765
$ws.Range("H125").Value = 1969.3334
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1969.3334
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17724.0006
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -22644.0006
$ws.Range("H127").Value = 831.38464
$ws.Range("I127").Value = 437.2857
$ws.Range("J127").Value = 1291.1666
$ws.Range("K127").Value = 1311.8571
$ws.Range("L127").Value = 3873.4998
$ws.Range("M127").Value = 3648.1429
$ws.Range("N127").Value = -13793.4998

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1702.75
$ws.Range("I97").Value = 457.77777
$ws.Range("J97").Value = 5437.6665
$ws.Range("K97").Value = 457.77777
$ws.Range("L97").Value = 5437.6665
$ws.Range("M97").Value = 38.22223000000002
$ws.Range("N97").Value = -6429.6665
$ws.Range("H102").Value = 5684018.5
$ws.Range("I102").Value = 6946619.5
$ws.Range("J102").Value = 2312.5
$ws.Range("K102").Value = 6946619.5
$ws.Range("L102").Value = 2312.5
$ws.Range("M102").Value = -6944997.5
$ws.Range("N102").Value = -5556.5
$ws.Range("H139").Value = 131300
$ws.Range("J139").Value = 131300
$ws.Range("L139").Value = 131300
$ws.Range("N139").Value = -141580

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 1480
$ws.Range("I97").Value = 1480
$ws.Range("K97").Value = 1480
$ws.Range("M97").Value = -489
$ws.Range("H99").Value = 1667.2667
$ws.Range("I99").Value = 1620.9
$ws.Range("J99").Value = 1760
$ws.Range("K99").Value = 1620.9
$ws.Range("L99").Value = 1760
$ws.Range("M99").Value = -122.9000000000001
$ws.Range("N99").Value = -4756

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 172.85715
$ws.Range("J86").Value = 300
$ws.Range("L86").Value = 900
$ws.Range("N86").Value = -3272
$ws.Range("H89").Value = 172.85715
$ws.Range("J89").Value = 300
$ws.Range("L89").Value = 2700
$ws.Range("N89").Value = -14556
$ws.Range("H121").Value = 1438.5555
$ws.Range("J121").Value = 1626.2667
$ws.Range("L121").Value = 4878.800099999999
$ws.Range("N121").Value = -7498.800099999999
$ws.Range("H131").Value = 1973184.9
$ws.Range("I131").Value = 514.2857
$ws.Range("J131").Value = 2404706.5
$ws.Range("K131").Value = 1542.8571
$ws.Range("L131").Value = 7214119.5
$ws.Range("M131").Value = 3497.1429
$ws.Range("N131").Value = -7224199.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2799.3684
$ws.Range("I100").Value = 1982.2858
$ws.Range("J100").Value = 2983.8708
$ws.Range("K100").Value = 1982.2858
$ws.Range("L100").Value = 2983.8708
$ws.Range("M100").Value = -1441.2858
$ws.Range("N100").Value = -4065.8708
$ws.Range("H122").Value = 7200.5557
$ws.Range("I122").Value = 7828.5713
$ws.Range("J122").Value = 5002.5
$ws.Range("K122").Value = 23485.7139
$ws.Range("L122").Value = 15007.5
$ws.Range("M122").Value = -21035.7139
$ws.Range("N122").Value = -19907.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1133333.4
$ws.Range("J2").Value = 1133333.4
$ws.Range("L2").Value = 1133333.4
$ws.Range("N2").Value = -1133557.4
$ws.Range("H11").Value = 25000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 25000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 25000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -25284
$ws.Range("H20").Value = 6055.5557
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 6750
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 6750
$ws.Range("M20").Value = -260
$ws.Range("N20").Value = -7230
$ws.Range("H22").Value = 60000
$ws.Range("J22").Value = 60000
$ws.Range("L22").Value = 60000
$ws.Range("N22").Value = -60586
$ws.Range("H30").Value = 11875
$ws.Range("J30").Value = 11875
$ws.Range("L30").Value = 11875
$ws.Range("N30").Value = -12089
$ws.Range("H33").Value = 70000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 70000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 70000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -70500
$ws.Range("H34").Value = 70000
$ws.Range("J34").Value = 70000
$ws.Range("L34").Value = 70000
$ws.Range("N34").Value = -70406
$ws.Range("H36").Value = 70000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 70000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 70000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -70500
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30298
$ws.Range("H44").Value = 43000
$ws.Range("J44").Value = 43000
$ws.Range("L44").Value = 43000
$ws.Range("N44").Value = -44108
$ws.Range("H47").Value = 36666.668
$ws.Range("J47").Value = 36666.668
$ws.Range("L47").Value = 36666.668
$ws.Range("N47").Value = -37810.668
$ws.Range("H48").Value = 70000
$ws.Range("J48").Value = 70000
$ws.Range("L48").Value = 70000
$ws.Range("N48").Value = -71138
$ws.Range("H49").Value = 30000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 30000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -30460
$ws.Range("H50").Value = 65000
$ws.Range("J50").Value = 65000
$ws.Range("L50").Value = 65000
$ws.Range("N50").Value = -66262
$ws.Range("H54").Value = 27538.5
$ws.Range("I54").Value = 9000
$ws.Range("J54").Value = 32173.125
$ws.Range("K54").Value = 9000
$ws.Range("L54").Value = 32173.125
$ws.Range("M54").Value = -8480
$ws.Range("N54").Value = -33213.125
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
